$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 70, shifting rows 70:114 down to 71:115.
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with fresh data (copy constant columns from the
# row that used to be at 70, now at 71; set the changed metrics).
$ws.Range("A70").Value = 9
$ws.Range("B70").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = 44673
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = 100112022
$ws.Range("G70").Value = "Arveja Verde"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 18
$ws.Range("K70").Value = 22000
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = 23000
$ws.Range("N70").Value = "$/saco 25 kilos"
$ws.Range("O70").Value = "Carahue"
$ws.Range("P70").Value = 920
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"
